$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2198.2942
$ws.Range("I15").Value = 2198.2942
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 6594.882599999999
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6425.882599999999
$ws.Range("H107").Value = 991.9167
$ws.Range("I107").Value = 991.9167
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 991.9167
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 928.0833
$ws.Range("H111").Value = 842.5
$ws.Range("I111").Value = 375
$ws.Range("J111").Value = 998.3333
$ws.Range("K111").Value = 1125
$ws.Range("L111").Value = 2994.9999
$ws.Range("M111").Value = 1942
$ws.Range("N111").Value = -9128.999899999999
$ws.Range("H131").Value = 13756
$ws.Range("I131").Value = 1155.4166
$ws.Range("J131").Value = 43997.4
$ws.Range("K131").Value = 3466.2498
$ws.Range("L131").Value = 131992.2
$ws.Range("M131").Value = 1573.7502
$ws.Range("N131").Value = -142072.2
$ws.Range("H137").Value = 26811.309
$ws.Range("I137").Value = 32091.088
$ws.Range("J137").Value = 4372.25
$ws.Range("K137").Value = 96273.264
$ws.Range("L137").Value = 13116.75
$ws.Range("M137").Value = -93723.264
$ws.Range("N137").Value = -18216.75
$ws.Range("H138").Value = 3016.3635
$ws.Range("I138").Value = 2974.9048
$ws.Range("J138").Value = 3041.9707
$ws.Range("K138").Value = 8924.714399999999
$ws.Range("L138").Value = 9125.9121
$ws.Range("M138").Value = -3784.714399999999
$ws.Range("N138").Value = -19405.9121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4292.6665
$ws.Range("I102").Value = 4115.8335
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 4115.8335
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -2493.8335
$ws.Range("N102").Value = -8244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 77763
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 77763
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 77763
$ws.Range("N81").Value = -79885
$ws.Range("H84").Value = 77763
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 77763
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 233289
$ws.Range("N84").Value = -243897

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2401
$ws.Range("I31").Value = 1543.4615
$ws.Range("J31").Value = 3993.5715
$ws.Range("K31").Value = 1543.4615
$ws.Range("L31").Value = 3993.5715
$ws.Range("M31").Value = -1248.4615
$ws.Range("N31").Value = -4583.5715
$ws.Range("H34").Value = 2401
$ws.Range("I34").Value = 1543.4615
$ws.Range("J34").Value = 3993.5715
$ws.Range("K34").Value = 1543.4615
$ws.Range("L34").Value = 3993.5715
$ws.Range("M34").Value = -1341.4615
$ws.Range("N34").Value = -4397.5715
$ws.Range("H99").Value = 2179
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2179
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 2179
$ws.Range("N99").Value = -5175
$ws.Range("M99").ClearContents()
$ws.Range("H103").Value = 26170.166
$ws.Range("I103").Value = 17404.4
$ws.Range("J103").Value = 69999
$ws.Range("K103").Value = 17404.4
$ws.Range("L103").Value = 69999
$ws.Range("M103").Value = -16232.4
$ws.Range("N103").Value = -72343
$ws.Range("H126").Value = 2179
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2179
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6537
$ws.Range("N126").Value = -11477
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 1721.579
$ws.Range("I132").Value = 1650.5555
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 4951.666499999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2421.666499999999
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 4349.5
$ws.Range("I34").Value = 699
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 2097
$ws.Range("L34").Value = 24000
$ws.Range("M34").Value = -2013
$ws.Range("N34").Value = -24168
$ws.Range("H39").Value = 9162.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 9162.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 27487.5
$ws.Range("N39").Value = -28075.5
$ws.Range("M39").ClearContents()
$ws.Range("H55").Value = 4346.3237
$ws.Range("I55").Value = 3056
$ws.Range("J55").Value = 4622.8213
$ws.Range("K55").Value = 9168
$ws.Range("L55").Value = 13868.4639
$ws.Range("M55").Value = -8991
$ws.Range("N55").Value = -14222.4639
$ws.Range("H107").Value = 711.2222
$ws.Range("I107").Value = 259.8
$ws.Range("J107").Value = 1275.5
$ws.Range("K107").Value = 779.4000000000001
$ws.Range("L107").Value = 3826.5
$ws.Range("M107").Value = 1140.6
$ws.Range("N107").Value = -7666.5
$ws.Range("H111").Value = 6331.3335
$ws.Range("I111").Value = 3997.5
$ws.Range("J111").Value = 10999
$ws.Range("K111").Value = 11992.5
$ws.Range("L111").Value = 32997
$ws.Range("M111").Value = -8925.5
$ws.Range("N111").Value = -39131
$ws.Range("H119").Value = 10499
$ws.Range("I119").Value = 999
$ws.Range("J119").Value = 19999
$ws.Range("K119").Value = 2997
$ws.Range("L119").Value = 59997
$ws.Range("M119").Value = 1841
$ws.Range("N119").Value = -69673
$ws.Range("H127").Value = 6827.6665
$ws.Range("I127").Value = 2500
$ws.Range("J127").Value = 8991.5
$ws.Range("K127").Value = 7500
$ws.Range("L127").Value = 26974.5
$ws.Range("M127").Value = -2540
$ws.Range("N127").Value = -36894.5
$ws.Range("H131").Value = 4778209
$ws.Range("I131").Value = 17873.166
$ws.Range("J131").Value = 6682343.5
$ws.Range("K131").Value = 53619.49800000001
$ws.Range("L131").Value = 20047030.5
$ws.Range("M131").Value = -48579.49800000001
$ws.Range("N131").Value = -20057110.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 7505
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 7505
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 7505
$ws.Range("N23").Value = -7951
$ws.Range("H44").Value = 28330.334
$ws.Range("I44").Value = 4999.5
$ws.Range("J44").Value = 39995.75
$ws.Range("K44").Value = 4999.5
$ws.Range("L44").Value = 39995.75
$ws.Range("M44").Value = -4403.5
$ws.Range("N44").Value = -41187.75
$ws.Range("H82").Value = 35998
$ws.Range("I82").Value = 23997.5
$ws.Range("J82").Value = 59999
$ws.Range("K82").Value = 23997.5
$ws.Range("L82").Value = 59999
$ws.Range("M82").Value = -23614.5
$ws.Range("N82").Value = -60765
$ws.Range("H85").Value = 35998
$ws.Range("I85").Value = 23997.5
$ws.Range("J85").Value = 59999
$ws.Range("K85").Value = 23997.5
$ws.Range("L85").Value = 59999
$ws.Range("M85").Value = -22671.5
$ws.Range("N85").Value = -62651
$ws.Range("H102").Value = 3261.4375
$ws.Range("I102").Value = 3106.4614
$ws.Range("J102").Value = 3933
$ws.Range("K102").Value = 3106.4614
$ws.Range("L102").Value = 3933
$ws.Range("M102").Value = -1484.4614
$ws.Range("N102").Value = -7177
$ws.Range("H132").Value = 22447.04
$ws.Range("I132").Value = 34084.906
$ws.Range("J132").Value = 3826.45
$ws.Range("K132").Value = 102254.718
$ws.Range("L132").Value = 11479.35
$ws.Range("M132").Value = -99724.71800000001
$ws.Range("N132").Value = -16539.35

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H55").Value = 975.4545000000001
$ws.Range("I55").Value = 376.8
$ws.Range("J55").Value = 1474.3334
$ws.Range("K55").Value = 376.8
$ws.Range("L55").Value = 1474.3334
$ws.Range("M55").Value = -203.8
$ws.Range("N55").Value = -1820.3334
$ws.Range("H132").Value = 63315.95
$ws.Range("I132").Value = 77795.94
$ws.Range("J132").Value = 5396
$ws.Range("K132").Value = 233387.82
$ws.Range("L132").Value = 16188
$ws.Range("M132").Value = -230857.82
$ws.Range("N132").Value = -21248
$ws.Range("H136").Value = 2309.5334
$ws.Range("I136").Value = 2180.45
$ws.Range("J136").Value = 2567.7
$ws.Range("K136").Value = 6541.349999999999
$ws.Range("L136").Value = 7703.099999999999
$ws.Range("M136").Value = -3991.349999999999
$ws.Range("N136").Value = -12803.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 17348
$ws.Range("I39").Value = 12022
$ws.Range("J39").Value = 28000
$ws.Range("K39").Value = 12022
$ws.Range("L39").Value = 28000
$ws.Range("M39").Value = -11609
$ws.Range("N39").Value = -28826
$ws.Range("H41").Value = 29097.2
$ws.Range("I41").Value = 28996
$ws.Range("J41").Value = 29249
$ws.Range("K41").Value = 28996
$ws.Range("L41").Value = 29249
$ws.Range("M41").Value = -28606
$ws.Range("N41").Value = -30029
$ws.Range("H122").Value = 5906
$ws.Range("I122").Value = 6513.3335
$ws.Range("J122").Value = 4995
$ws.Range("K122").Value = 19540.0005
$ws.Range("L122").Value = 14985
$ws.Range("M122").Value = -17090.0005
$ws.Range("N122").Value = -19885

